# The commit swaps the two theme parts of the deck:
#   ppt/theme/theme1.xml  "Office Theme"  <->  ppt/theme/theme2.xml  "Integral"
# theme2.xml is the theme actually in force for the slides/slide master
# (ppt/slideMasters/slideMaster1.xml and the presentation's own theme
# relationship both point at theme2.xml), so the externally-visible effect
# of the swap is that the deck's colour palette changes from the green/
# yellow "Integral" scheme to the default blue/orange "Office" scheme.
#
# Drive that through the real PowerPoint theme-colours API: each of the
# twelve theme colour slots (dk1, lt1, dk2, lt2, accent1-6, hlink,
# folHlink) is exposed as ThemeColorScheme.Colors(1..12).RGB on the
# slide master's Theme - this is the supported, persisted way to edit a
# deck's colour scheme via COM automation.

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$colorScheme = $master.Theme.ThemeColorScheme

# Target palette: the stock "Office Theme" colours (what theme2.xml's
# <a:clrScheme> should contain after the swap), in the fixed COM order
# dk1, lt1, dk2, lt2, accent1, accent2, accent3, accent4, accent5,
# accent6, hlink, folHlink. Values are packed as 0x00BBGGRR, matching
# the PowerPoint COM RGB color convention (long integer form of RGB()).
$officeThemeColors = @(
    0x000000,   # dk1       000000
    0xFFFFFF,   # lt1       FFFFFF
    0x6A5444,   # dk2       44546A
    0xE6E6E7,   # lt2       E7E6E6
    0xD59B5B,   # accent1   5B9BD5
    0x317DED,   # accent2   ED7D31
    0xA5A5A5,   # accent3   A5A5A5
    0x00C0FF,   # accent4   FFC000
    0xC47244,   # accent5   4472C4
    0x47AD70,   # accent6   70AD47
    0xC16305,   # hlink     0563C1
    0x724F95    # folHlink  954F72
)

for ($i = 1; $i -le $officeThemeColors.Count; $i++) {
    $colorScheme.Colors($i).RGB = $officeThemeColors[$i - 1]
}
